$wb = $excel.ActiveWorkbook

# --- Sheet "Relations" (1st sheet) ---------------------------------------
# "Park" the mistake rows: move A6:C8 content to D6:F8 (same values), and
# clear out the old A6:C8 area.
$wsRelations = $wb.Worksheets.Item(1)

$wsRelations.Range("D6").Value = $wsRelations.Range("A6").Value2

$wsRelations.Range("D7").Value = $wsRelations.Range("A7").Value2
$wsRelations.Range("E7").Value = $wsRelations.Range("B7").Value2
$wsRelations.Range("F7").Value = $wsRelations.Range("C7").Value2

$wsRelations.Range("D8").Value = $wsRelations.Range("A8").Value2
$wsRelations.Range("E8").Value = $wsRelations.Range("B8").Value2
$wsRelations.Range("F8").Value = $wsRelations.Range("C8").Value2

$wsRelations.Range("A6:C8").Clear() | Out-Null

$wsRelations.Range("D6:F8").Select() | Out-Null

# --- Sheet "Rules" (2nd sheet) --------------------------------------------
# Park the old "undeclared = w" mistake row down to row 5 (columns D:F) and
# introduce the new "v = w" rule on row 4 (columns D:F); clear old A4:C4.
$wsRules = $wb.Worksheets.Item(2)

$wsRules.Range("D5").Value = $wsRules.Range("A4").Value2
$wsRules.Range("E5").Value = $wsRules.Range("B4").Value2
$wsRules.Range("F5").Value = $wsRules.Range("C4").Value2

$wsRules.Range("D4").Value = "v = w"
$wsRules.Range("E4").Value = "v"
$wsRules.Range("F4").Value = "w"

$wsRules.Range("A4:C4").Clear() | Out-Null

# --- Switch the active tab back to "Relations" ----------------------------
$wsRelations.Activate() | Out-Null
